$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updated s_val data
$ws.Range("B2").Value = 3.182878228561681
$ws.Range("C2").Value = 1.65323645889881
$ws.Range("D2").Value = 0.1529057820181812
$ws.Range("E2").Value = 0.4998867070740569
$ws.Range("G2").Value = 5.488907176552729

# Row 3 updated s_val data
$ws.Range("B3").Value = 1.505614041169197
$ws.Range("C3").Value = 10990084.13351303
$ws.Range("D3").Value = 16.98373111632243
$ws.Range("E3").Value = 6.48142807727062
$ws.Range("G3").Value = 10990109.10428627
